# Auto-optimize exam scheduling: dynamically adjusts exams per slot (1-4)
# to guarantee all courses are scheduled within date range.
# This updates the Section_A / Section_B timetables and the
# Elective_Coordination slot times to reflect the re-optimized schedule.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")
$wsElective = $wb.Worksheets.Item("Elective_Coordination")

# ---- Section_A timetable updates ----
$wsA.Range("B2").Value = "MA161"
$wsA.Range("C2").Value = "MA161"
$wsA.Range("D2").Value = "MA162"
$wsA.Range("E2").Value = "C202"
$wsA.Range("F2").Value = "DS161"

$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "DS161"
$wsA.Range("E3").Value = "CS151 (Elective)"
$wsA.Range("F3").Value = "CS161"

$wsA.Range("B5").Value = "C202"
$wsA.Range("C5").Value = "CS161"
$wsA.Range("D5").Value = "Free"
$wsA.Range("F5").Value = "MA162"

$wsA.Range("B6").Value = "Free"

$wsA.Range("B7").Value = "EC161"
$wsA.Range("C7").Value = "EC161"
$wsA.Range("D7").Value = "C202"
$wsA.Range("E7").Value = "CS161"
$wsA.Range("F7").Value = "CS151 (Elective)"

$wsA.Range("C8").Value = "CS151 (Tutorial)"

# ---- Section_B timetable updates ----
$wsB.Range("B2").Value = "MA162"
$wsB.Range("C2").Value = "C202"
$wsB.Range("E2").Value = "DS161"
$wsB.Range("F2").Value = "CS161"

$wsB.Range("B3").Value = "MA161"
$wsB.Range("C3").Value = "DS161"
$wsB.Range("E3").Value = "CS151 (Elective)"
$wsB.Range("F3").Value = "Free"

$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "MA162"
$wsB.Range("D5").Value = "CS161"
$wsB.Range("E5").Value = "EC161"

$wsB.Range("B6").Value = "Free"

$wsB.Range("B7").Value = "EC161"
$wsB.Range("C7").Value = "CS161"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "MA161"
$wsB.Range("F7").Value = "CS151 (Elective)"

$wsB.Range("C8").Value = "CS151 (Tutorial)"

# ---- Elective_Coordination slot-time updates ----
$wsElective.Range("D10").Value = "10:30-12:00"
$wsElective.Range("D11").Value = "15:30-17:00"
$wsElective.Range("C12").Value = "Tue"
$wsElective.Range("D12").Value = "17:00-18:00"
